$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows (in column H) whose value changes from "NULL" to "yes"
$rows = @(3,6,7,8,9,10,11,12,13,14,15,16,17,20,21,22,23,24,27,28,29,30,31,32,33,34,36,37,38,43,44,45,48,49,50,51,52,54,55,56)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 8).Value = "yes"
}

# Rows 44 and 45 also pick up the same formatting style as the cells above them (row 43)
$ws.Range("H43").Copy()
$ws.Range("H44:H45").PasteSpecial(-4122)  # xlPasteFormats

# Update the selection / scroll position to match the saved view
$ws.Range("I54").Select()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
